$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Contribuição" column (E) holds percentages stored as plain text
# (e.g. "43%"), not numeric percentage values. Force the number format to
# text ("@") before assigning so Excel doesn't coerce the string into a
# numeric percentage value.
$pctCells = @("E2","E3","E5","E6","E8","E9","E10","E11")
foreach ($cell in $pctCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Row 2: K. Jorge, Cruzeiro
$ws.Range("D2").Value = 1330
$ws.Range("E2").Value = "43%"

# Row 3: G. d. Arrascaeta, Flamengo
$ws.Range("D3").Value = 1031
$ws.Range("E3").Value = "32%"

# Row 4: P. Vegetti, V. d. Gama (E unchanged)
$ws.Range("D4").Value = 918

# Row 5: Reinaldo, Mirassol
$ws.Range("D5").Value = 727
$ws.Range("E5").Value = "26%"

# Row 6: now P. Raul / Ceara (was M. Braithwaite / Gremio)
$ws.Range("B6").Value = "P. Raul"
$ws.Range("C6").Value = "Ceara"
$ws.Range("D6").Value = 717
$ws.Range("E6").Value = "41%"

# Row 7: now M. Braithwaite / Gremio (was P. Raul / Ceara), D7/E7 unchanged
$ws.Range("B7").Value = "M. Braithwaite"
$ws.Range("C7").Value = "Gremio"

# Row 8: A. Silva, S. Paulo
$ws.Range("D8").Value = 520
$ws.Range("E8").Value = "25%"

# Row 9: Y. Alberto, Corinthians
$ws.Range("D9").Value = 517
$ws.Range("E9").Value = "29%"

# Row 10: now R. Kayzer / Vitoria (was Pedro / Flamengo)
$ws.Range("B10").Value = "R. Kayzer"
$ws.Range("C10").Value = "Vitoria"
$ws.Range("D10").Value = 516
$ws.Range("E10").Value = "31%"

# Row 11: now Pedro / Flamengo (was I. Pitta / R. B. Bragantino)
$ws.Range("B11").Value = "Pedro"
$ws.Range("C11").Value = "Flamengo"
$ws.Range("D11").Value = 531
$ws.Range("E11").Value = "16%"
